# fix(working-calendar): crud working calendar
#
# Rework the header row (D1/E1 renamed, F1/G1 turned into blank spare
# columns) and append three working-calendar rows (working day / event /
# holiday) beneath it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- header row -------------------------------------------------------
# A1:C1 (comp_code / branch_code / date) are unchanged; D1/E1 get new
# labels and F1/G1 become blank placeholder columns.
$ws.Range("D1").Value = "type"
$ws.Range("E1").Value = "description"

# --- row 2: working day ------------------------------------------------
# The date column holds date-LOOKING text ("2024-01-01"), not a real
# date. A leading apostrophe forces text entry so Excel doesn't coerce
# it into a date serial number.
$ws.Range("A2").Value = "C001"
$ws.Range("B2").Value = "B001"
$ws.Range("C2").Value = "'2024-01-01"
$ws.Range("D2").Value = "working day"
$ws.Range("E2").Value = "hari kerja"

# --- row 3: event --------------------------------------------------------
$ws.Range("A3").Value = "C001"
$ws.Range("B3").Value = "B001"
$ws.Range("C3").Value = "'2024-01-02"
$ws.Range("D3").Value = "event"
$ws.Range("E3").Value = "event"

# --- row 4: holiday -------------------------------------------------------
$ws.Range("A4").Value = "C001"
$ws.Range("B4").Value = "B001"
$ws.Range("C4").Value = "'2024-01-03"
$ws.Range("D4").Value = "holiday"
$ws.Range("E4").Value = "libur"

# F1:G4 stay present but blank (spare columns). A leading apostrophe
# forces an explicit, empty text entry instead of clearing the cell
# outright.
$ws.Range("F1:G4").Value = "'"

# Resetting the style afterwards drops the quote-prefix flag that
# apostrophe text-entry leaves behind, so the cells end up styled
# exactly like freshly-typed text (no visible "number stored as text"
# marker, matching the rest of the sheet).
$ws.Range("C2:C4").Style = "Normal"
$ws.Range("F1:G4").Style = "Normal"
